$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header label for column E
$ws.Range("E1").Value = "strength (RMS)"

# Updated measurement values (columns B, C, D, E) per row, rows 2-19
$data = @(
    @{ Row = 2;  B = 16;    C = 27;   D = 11;   E = 41 }
    @{ Row = 3;  B = 15.6;  C = 30;   D = 14.4; E = 44.6 }
    @{ Row = 4;  B = 18.8;  C = 30.4; D = 11.6; E = 42.4 }
    @{ Row = 5;  B = 14.8;  C = 30;   D = 15.2; E = 48.2 }
    @{ Row = 6;  B = 15.6;  C = 30.4; D = 14.8; E = 40 }
    @{ Row = 7;  B = 15.2;  C = 29.6; D = 14.4; E = 50.2 }
    @{ Row = 8;  B = 14.8;  C = 30;   D = 15.2; E = 48.4 }
    @{ Row = 9;  B = 16;    C = 30.4; D = 14.4; E = 44.2 }
    @{ Row = 10; B = 16;    C = 30;   D = 14;   E = 44.6 }
    @{ Row = 11; B = 14.4;  C = 30;   D = 15.6; E = 47.2 }
    @{ Row = 12; B = 16;    C = 30.4; D = 14.4; E = 44 }
    @{ Row = 13; B = 16;    C = 30;   D = 14;   E = 46.5 }
    @{ Row = 14; B = 16;    C = 30.4; D = 14.4; E = 41.2 }
    @{ Row = 15; B = 15.6;  C = 30.4; D = 14.8; E = 41.2 }
    @{ Row = 16; B = 16;    C = 30;   D = 14;   E = 48.75 }
    @{ Row = 17; B = 15.5;  C = 30;   D = 14.5; E = 46.75 }
    @{ Row = 18; B = 14.8;  C = 30;   D = 15.2; E = 44.6 }
    @{ Row = 19; B = 16;    C = 30;   D = 14;   E = 43.5 }
)

foreach ($entry in $data) {
    $r = $entry.Row
    $ws.Cells.Item($r, 2).Value = $entry.B
    $ws.Cells.Item($r, 3).Value = $entry.C
    $ws.Cells.Item($r, 4).Value = $entry.D
    $ws.Cells.Item($r, 5).Value = $entry.E
}
